$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 16.89431068226313
$ws.Range("C2").Value = 11.27777840111306
$ws.Range("D2").Value = 5.350536554766854
$ws.Range("F2").Value = 29.07368784240456
$ws.Range("G2").Value = 37.50544705067728
$ws.Range("H2").Value = 16.07299276724145
$ws.Range("L2").Value = 10.71929615084124
$ws.Range("M2").Value = 15.95783662952572
$ws.Range("N2").Value = 18.75475445371866
# Row 3
$ws.Range("B3").Value = 16.43596285531251
$ws.Range("C3").Value = 11.00974897504455
$ws.Range("D3").Value = 5.358677960774445
$ws.Range("F3").Value = 28.94556821552616
$ws.Range("G3").Value = 37.22429545868221
$ws.Range("H3").Value = 16.09118361987715
$ws.Range("L3").Value = 10.73452150675075
$ws.Range("M3").Value = 15.88106756311351
$ws.Range("N3").Value = 18.82637848821093
# Row 4
$ws.Range("B4").Value = 16.15203426945391
$ws.Range("C4").Value = 10.83998804464156
$ws.Range("D4").Value = 5.363924950982267
$ws.Range("F4").Value = 28.87602119513166
$ws.Range("G4").Value = 37.06476107081937
$ws.Range("H4").Value = 16.10608002602714
$ws.Range("L4").Value = 10.74540066189807
$ws.Range("M4").Value = 15.83691825520078
$ws.Range("N4").Value = 18.87229701227112
# Row 5
$ws.Range("B5").Value = 16.03589088573713
$ws.Range("C5").Value = 10.76955867769252
$ws.Range("D5").Value = 5.366125802953968
$ws.Range("F5").Value = 28.84999319093487
$ws.Range("G5").Value = 37.00310350712443
$ws.Range("H5").Value = 16.11308546372733
$ws.Range("L5").Value = 10.750218992836
$ws.Range("M5").Value = 15.81969199735325
$ws.Range("N5").Value = 18.89149896006376
# Row 6
$ws.Range("B6").Value = 16.0165844766675
$ws.Range("C6").Value = 10.75779007984713
$ws.Range("D6").Value = 5.366495045155007
$ws.Range("F6").Value = 28.84581146841728
$ws.Range("G6").Value = 36.99306953796672
$ws.Range("H6").Value = 16.11430511404281
$ws.Range("L6").Value = 10.75104232651552
$ws.Range("M6").Value = 15.81687817397154
$ws.Range("N6").Value = 18.89471705939838
# Row 7
$ws.Range("B7").Value = 16.15046944522045
$ws.Range("C7").Value = 10.83904319760132
$ws.Range("D7").Value = 5.363954378391489
$ws.Range("F7").Value = 28.8756607837289
$ws.Range("G7").Value = 37.0639158796069
$ws.Range("H7").Value = 16.10617072110946
$ws.Range("L7").Value = 10.74546408473477
$ws.Range("M7").Value = 15.83668282120075
$ws.Range("N7").Value = 18.87255399117566
# Row 8
$ws.Range("B8").Value = 16.73690174824507
$ws.Range("C8").Value = 11.18647061913488
$ws.Range("D8").Value = 5.353292370552574
$ws.Range("F8").Value = 29.02763263468865
$ws.Range("G8").Value = 37.40582747534734
$ws.Range("H8").Value = 16.07849018629358
$ws.Range("L8").Value = 10.72422820808953
$ws.Range("M8").Value = 15.93075523120465
$ws.Range("N8").Value = 18.77904856503973
# Row 9
$ws.Range("B9").Value = 17.85926327265036
$ws.Range("C9").Value = 11.82438643190792
$ws.Range("D9").Value = 5.334341090591135
$ws.Range("F9").Value = 29.39695517248684
$ws.Range("G9").Value = 38.17715651384776
$ws.Range("H9").Value = 16.05386635257521
$ws.Range("L9").Value = 10.69472618635691
$ws.Range("M9").Value = 16.13828280215571
$ws.Range("N9").Value = 18.61100830596449
# Row 10
$ws.Range("B10").Value = 18.65737777383196
$ws.Range("C10").Value = 12.26391799829775
$ws.Range("D10").Value = 5.321594072783778
$ws.Range("F10").Value = 29.71012307276252
$ws.Range("G10").Value = 38.80078079360953
$ws.Range("H10").Value = 16.05395037052321
$ws.Range("L10").Value = 10.68044581113862
$ws.Range("M10").Value = 16.30387489708052
$ws.Range("N10").Value = 18.49677954647522
# Row 11
$ws.Range("B11").Value = 19.01290860374274
$ws.Range("C11").Value = 12.45704643364208
$ws.Range("D11").Value = 5.316047135324506
$ws.Range("F11").Value = 29.86125442464778
$ws.Range("G11").Value = 39.0957605116188
$ws.Range("H11").Value = 16.05794758881534
$ws.Range("L11").Value = 10.6755525266439
$ws.Range("M11").Value = 16.38184066741271
$ws.Range("N11").Value = 18.44679437081693
# Row 12
$ws.Range("B12").Value = 19.14632173947444
$ws.Range("C12").Value = 12.52916191860088
$ws.Range("D12").Value = 5.313982594833497
$ws.Range("F12").Value = 29.91969081166974
$ws.Range("G12").Value = 39.20898762828861
$ws.Range("H12").Value = 16.06003080844652
$ws.Range("L12").Value = 10.6739297448487
$ws.Range("M12").Value = 16.41172420267452
$ws.Range("N12").Value = 18.42814895122515
# Row 13
$ws.Range("B13").Value = 19.11764505678972
$ws.Range("C13").Value = 12.5136764037098
$ws.Range("D13").Value = 5.31442563475488
$ws.Range("F13").Value = 29.90705251551515
$ws.Range("G13").Value = 39.1845359395688
$ws.Range("H13").Value = 16.05955681887253
$ws.Range("L13").Value = 10.67426900720903
$ws.Range("M13").Value = 16.4052725745428
$ws.Range("N13").Value = 18.43215201945996
# Row 14
$ws.Range("B14").Value = 19.02390973748179
$ws.Range("C14").Value = 12.46300002616117
$ws.Range("D14").Value = 5.31587656489514
$ws.Range("F14").Value = 29.86603808809383
$ws.Range("G14").Value = 39.10504579512594
$ws.Range("H14").Value = 16.05810756174575
$ws.Range("L14").Value = 10.67541440781461
$ws.Range("M14").Value = 16.38429211163483
$ws.Range("N14").Value = 18.44525474046788
# Row 15
$ws.Range("B15").Value = 18.96633163821321
$ws.Range("C15").Value = 12.43182563487907
$ws.Range("D15").Value = 5.316769978181191
$ws.Range("F15").Value = 29.84107140640515
$ws.Range("G15").Value = 39.05655135779361
$ws.Range("H15").Value = 16.05729402386944
$ws.Range("L15").Value = 10.67614596750722
$ws.Range("M15").Value = 16.37148720895241
$ws.Range("N15").Value = 18.45331732877638
# Row 16
$ws.Range("B16").Value = 18.63397954064667
$ws.Range("C16").Value = 12.25115612024175
$ws.Range("D16").Value = 5.321961623643358
$ws.Range("F16").Value = 29.70041734253446
$ws.Range("G16").Value = 38.78172199369649
$ws.Range("H16").Value = 16.05376884006379
$ws.Range("L16").Value = 10.68079782950392
$ws.Range("M16").Value = 16.29883105252893
$ws.Range("N16").Value = 18.50008586523148
# Row 17
$ws.Range("B17").Value = 18.42806230596655
$ws.Range("C17").Value = 12.13854730862042
$ws.Range("D17").Value = 5.325210838386067
$ws.Range("F17").Value = 29.61632322580908
$ws.Range("G17").Value = 38.61594556368851
$ws.Range("H17").Value = 16.05262056874009
$ws.Range("L17").Value = 10.68406194997012
$ws.Range("M17").Value = 16.25492051586976
$ws.Range("N17").Value = 18.52928236257947
# Row 18
$ws.Range("B18").Value = 18.30892218241446
$ws.Range("C18").Value = 12.07313812172143
$ws.Range("D18").Value = 5.327103409974261
$ws.Range("F18").Value = 29.56877315088757
$ws.Range("G18").Value = 38.52166416868589
$ws.Range("H18").Value = 16.05233276414727
$ws.Range("L18").Value = 10.68609028278447
$ws.Range("M18").Value = 16.22991376871779
$ws.Range("N18").Value = 18.54626169261854
# Row 19
$ws.Range("B19").Value = 18.26846713380376
$ws.Range("C19").Value = 12.0508830554967
$ws.Range("D19").Value = 5.327748281563472
$ws.Range("F19").Value = 29.5528153026996
$ws.Range("G19").Value = 38.48992863101428
$ws.Range("H19").Value = 16.0522993067287
$ws.Range("L19").Value = 10.68680296632396
$ws.Range("M19").Value = 16.22149034170859
$ws.Range("N19").Value = 18.55204264302973
# Row 20
$ws.Range("B20").Value = 18.45005621605678
$ws.Range("C20").Value = 12.15060118749931
$ws.Range("D20").Value = 5.324862501825708
$ws.Range("F20").Value = 29.62519074676553
$ws.Range("G20").Value = 38.6334828328859
$ws.Range("H20").Value = 16.05270423294604
$ws.Range("L20").Value = 10.6836988634833
$ws.Range("M20").Value = 16.25956919233683
$ws.Range("N20").Value = 18.52615507792778
# Row 21
$ws.Range("B21").Value = 19.05147616394963
$ws.Range("C21").Value = 12.47791282132938
$ws.Range("D21").Value = 5.315449417340214
$ws.Range("F21").Value = 29.8780526187337
$ws.Range("G21").Value = 39.12835340135909
$ws.Range("H21").Value = 16.05851778603708
$ws.Range("L21").Value = 10.67507173105181
$ws.Range("M21").Value = 16.39044498232438
$ws.Range("N21").Value = 18.44139848876112
# Row 22
$ws.Range("B22").Value = 19.43738131153346
$ws.Range("C22").Value = 12.68588065459866
$ws.Range("D22").Value = 5.309506943681065
$ws.Range("F22").Value = 30.05032268083294
$ws.Range("G22").Value = 39.4606190270787
$ws.Range("H22").Value = 16.06563692856879
$ws.Range("L22").Value = 10.67077499460755
$ws.Range("M22").Value = 16.47806633935732
$ws.Range("N22").Value = 18.38765329069771
# Row 23
$ws.Range("B23").Value = 19.23211290522619
$ws.Range("C23").Value = 12.57544029621202
$ws.Range("D23").Value = 5.312659459446978
$ws.Range("F23").Value = 29.95775151391029
$ws.Range("G23").Value = 39.28250729211093
$ws.Range("H23").Value = 16.06153358437092
$ws.Range("L23").Value = 10.67294560256851
$ws.Range("M23").Value = 16.43111689280957
$ws.Range("N23").Value = 18.4161878236093
# Row 24
$ws.Range("B24").Value = 18.44011511228294
$ws.Range("C24").Value = 12.14515370948217
$ws.Range("D24").Value = 5.325019908250444
$ws.Range("F24").Value = 29.62117925555872
$ws.Range("G24").Value = 38.62555103097174
$ws.Range("H24").Value = 16.05266524844977
$ws.Range("L24").Value = 10.68386254218565
$ws.Range("M24").Value = 16.25746678282677
$ws.Range("N24").Value = 18.52756831923054
# Row 25
$ws.Range("B25").Value = 17.5596337322939
$ws.Range("C25").Value = 11.65675071093564
$ws.Range("D25").Value = 5.339260136232431
$ws.Range("F25").Value = 29.28956607227584
$ws.Range("G25").Value = 37.95814184137536
$ws.Range("H25").Value = 16.05734134445677
$ws.Range("L25").Value = 10.70140794729809
$ws.Range("M25").Value = 16.07976594757442
$ws.Range("N25").Value = 18.65483868667901
